$d = $word.ActiveDocument

# Anchor: the paragraph ending in "...una pantalla de pausa." is currently the last
# paragraph in the document body.
$anchor = $d.Paragraphs.Last

function New-FollowingParagraph($prevRange) {
    $prevRange.InsertParagraphAfter()
    return $d.Paragraphs.Last
}

# ---------------------------------------------------------------------------
# 1-4: four blank paragraphs (same run formatting as the paragraph above: Arial,
#      sz 28 / szCs 28, firstLine indent 709 - inherited automatically).
# ---------------------------------------------------------------------------
$p = New-FollowingParagraph $anchor.Range
$p = New-FollowingParagraph $p.Range
$p = New-FollowingParagraph $p.Range
$p = New-FollowingParagraph $p.Range

# ---------------------------------------------------------------------------
# 5: "Día seis:" heading paragraph - bold/italic/underline, sz 32, single line
#     spacing (line=240, auto).
# ---------------------------------------------------------------------------
$p = New-FollowingParagraph $p.Range
$p.Format.LineSpacingRule = 0
$p.Format.LineSpacing = 12
$r = $p.Range
$r.InsertAfter("Día seis:")
$r = $d.Paragraphs.Last.Range
$r.Font.Name = "Arial"
$r.Font.NameAscii = "Arial"
$r.Font.NameBi = "Arial"
$r.Font.Bold = 1
$r.Font.BoldBi = 1
$r.Font.Italic = 1
$r.Font.ItalicBi = 1
$r.Font.Underline = 1
$r.Font.Size = 16
$r.Font.SizeBi = 16

# ---------------------------------------------------------------------------
# 6: "En realidad, no estoy seguro..." paragraph.
# ---------------------------------------------------------------------------
$p = New-FollowingParagraph $p.Range
$p.Format.LineSpacingRule = 0
$p.Format.LineSpacing = 12
$r = $p.Range
$r.InsertAfter("En realidad, no estoy seguro del día qué sea, pero cometí un grave error, intenté separar el script central del jugador en varios scripts que se encargasen de distintas funciones; SALIÓ MAL.")
$r = $d.Paragraphs.Last.Range
$r.Font.Name = "Arial"
$r.Font.NameAscii = "Arial"
$r.Font.NameBi = "Arial"
$r.Font.Bold = 0
$r.Font.BoldBi = 0
$r.Font.Italic = 0
$r.Font.ItalicBi = 0
$r.Font.Underline = 0
$r.Font.Size = 14
$r.Font.SizeBi = 14

# ---------------------------------------------------------------------------
# 7: "Tuve que borrar el proyecto, y clonarlo de GitHub, perdí algunos avances."
# ---------------------------------------------------------------------------
$p = New-FollowingParagraph $p.Range
$p.Format.LineSpacingRule = 0
$p.Format.LineSpacing = 12
$r = $p.Range
$r.InsertAfter("Tuve que borrar el proyecto, y clonarlo de GitHub, perdí algunos avances.")
$r = $d.Paragraphs.Last.Range
$r.Font.Name = "Arial"
$r.Font.NameAscii = "Arial"
$r.Font.NameBi = "Arial"
$r.Font.Bold = 0
$r.Font.BoldBi = 0
$r.Font.Italic = 0
$r.Font.ItalicBi = 0
$r.Font.Underline = 0
$r.Font.Size = 14
$r.Font.SizeBi = 14

# ---------------------------------------------------------------------------
# 8: "Lo pude salvar. Afortunadamente..."
# ---------------------------------------------------------------------------
$p = New-FollowingParagraph $p.Range
$p.Format.LineSpacingRule = 0
$p.Format.LineSpacing = 12
$r = $p.Range
$r.InsertAfter("Lo pude salvar. Afortunadamente, no habían muchas diferencias entre la versión que tuve que eliminar y la versión más reciente de GitHub: solo voy a tener que rehacer el sistema de pausa, y eso es todo.")
$r = $d.Paragraphs.Last.Range
$r.Font.Name = "Arial"
$r.Font.NameAscii = "Arial"
$r.Font.NameBi = "Arial"
$r.Font.Bold = 0
$r.Font.BoldBi = 0
$r.Font.Italic = 0
$r.Font.ItalicBi = 0
$r.Font.Underline = 0
$r.Font.Size = 14
$r.Font.SizeBi = 14

# ---------------------------------------------------------------------------
# 9: "Se rehicieron los menús de muerte y de pausa."
# ---------------------------------------------------------------------------
$p = New-FollowingParagraph $p.Range
$p.Format.LineSpacingRule = 0
$p.Format.LineSpacing = 12
$r = $p.Range
$r.InsertAfter("Se rehicieron los menús de muerte y de pausa.")
$r = $d.Paragraphs.Last.Range
$r.Font.Name = "Arial"
$r.Font.NameAscii = "Arial"
$r.Font.NameBi = "Arial"
$r.Font.Bold = 0
$r.Font.BoldBi = 0
$r.Font.Italic = 0
$r.Font.ItalicBi = 0
$r.Font.Underline = 0
$r.Font.Size = 14
$r.Font.SizeBi = 14

# ---------------------------------------------------------------------------
# 10: "Se les bajó la velocidad a los enemigos (2 -> 1)."
# ---------------------------------------------------------------------------
$p = New-FollowingParagraph $p.Range
$p.Format.LineSpacingRule = 0
$p.Format.LineSpacing = 12
$r = $p.Range
$r.InsertAfter("Se les bajó la velocidad a los enemigos (2 -> 1).")
$r = $d.Paragraphs.Last.Range
$r.Font.Name = "Arial"
$r.Font.NameAscii = "Arial"
$r.Font.NameBi = "Arial"
$r.Font.Bold = 0
$r.Font.BoldBi = 0
$r.Font.Italic = 0
$r.Font.ItalicBi = 0
$r.Font.Underline = 0
$r.Font.Size = 14
$r.Font.SizeBi = 14

# ---------------------------------------------------------------------------
# 11: "Se le subió la resistencia al jugador (100 -> 250)."
# ---------------------------------------------------------------------------
$p = New-FollowingParagraph $p.Range
$p.Format.LineSpacingRule = 0
$p.Format.LineSpacing = 12
$r = $p.Range
$r.InsertAfter("Se le subió la resistencia al jugador (100 -> 250).")
$r = $d.Paragraphs.Last.Range
$r.Font.Name = "Arial"
$r.Font.NameAscii = "Arial"
$r.Font.NameBi = "Arial"
$r.Font.Bold = 0
$r.Font.BoldBi = 0
$r.Font.Italic = 0
$r.Font.ItalicBi = 0
$r.Font.Underline = 0
$r.Font.Size = 14
$r.Font.SizeBi = 14

# ---------------------------------------------------------------------------
# 12: "Lo dejare por acá."
# ---------------------------------------------------------------------------
$p = New-FollowingParagraph $p.Range
$p.Format.LineSpacingRule = 0
$p.Format.LineSpacing = 12
$r = $p.Range
$r.InsertAfter("Lo dejare por acá.")
$r = $d.Paragraphs.Last.Range
$r.Font.Name = "Arial"
$r.Font.NameAscii = "Arial"
$r.Font.NameBi = "Arial"
$r.Font.Bold = 0
$r.Font.BoldBi = 0
$r.Font.Italic = 0
$r.Font.ItalicBi = 0
$r.Font.Underline = 0
$r.Font.Size = 14
$r.Font.SizeBi = 14

Write-Output ("Done. ParaCount=" + $d.Paragraphs.Count)
